$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the account.move test data rows (A:id, C:date/ref label, D:ref, E:journal_id, F:fiscalyear_id)
# Row 2
$ws.Range("A2").Value = "z0bug.move_01"
$ws.Range("C2").Value = "<002-12-99"
$ws.Range("D2").Value = "Capitale sociale"
$ws.Range("E2").Value = "external.MISC"
$ws.Range("F2").Value = "z0bug.fy_%(year)s"

# Row 3
$ws.Range("A3").Value = "z0bug.move_02"
$ws.Range("C3").Value = "<001-01-15"
$ws.Range("D3").Value = "Pagata RA"
$ws.Range("E3").Value = "external.MISC"
$ws.Range("F3").Value = "z0bug.fy_%(year)s"

# Row 4
$ws.Range("A4").Value = "z0bug.move_03"
$ws.Range("C4").Value = "<001-01-31"
$ws.Range("D4").Value = "Affitto 1.mo Trimestre 2020"
$ws.Range("E4").Value = "external.MISC"
$ws.Range("F4").Value = "z0bug.fy_%(year)s"

# Row 5
$ws.Range("A5").Value = "z0bug.move_04"
$ws.Range("C5").Value = "<001-04-05"
$ws.Range("D5").Value = "Affitto 2.do Trimestre 2020"
$ws.Range("E5").Value = "external.MISC"
$ws.Range("F5").Value = "z0bug.fy_%(year)s"

# Row 6
$ws.Range("A6").Value = "z0bug.move_05"
$ws.Range("C6").Value = "<001-07-05"
$ws.Range("D6").Value = "Affitto 3.zo Trimestre 2020"
$ws.Range("E6").Value = "external.MISC"
$ws.Range("F6").Value = "z0bug.fy_%(year)s"

# Row 7
$ws.Range("A7").Value = "z0bug.move_06"
$ws.Range("C7").Value = "<001-10-05"
$ws.Range("D7").Value = "Affitto 4.to Trimestre 2020"
$ws.Range("E7").Value = "external.MISC"
$ws.Range("F7").Value = "z0bug.fy_%(year)s"

# Row 8
$ws.Range("A8").Value = "z0bug.move_07"
$ws.Range("C8").Value = "####-01-31"
$ws.Range("D8").Value = "Affitto 1.mo Trimestre 2021"
$ws.Range("E8").Value = "external.MISC"
$ws.Range("F8").Value = "z0bug.fy_%(year)s"

# Row 9
$ws.Range("A9").Value = "z0bug.move_08"
$ws.Range("C9").Value = "####-04-05"
$ws.Range("D9").Value = "Affitto 2.do Trimestre 2021"
$ws.Range("E9").Value = "external.MISC"
$ws.Range("F9").Value = "z0bug.fy_%(year)s"

# Row 10
$ws.Range("A10").Value = "z0bug.move_09"
$ws.Range("C10").Value = "####-07-05"
$ws.Range("D10").Value = "Affitto 3.zo Trimestre 2021"
$ws.Range("E10").Value = "external.MISC"
$ws.Range("F10").Value = "z0bug.fy_%(year)s"

# Row 11
$ws.Range("A11").Value = "z0bug.move_10"
$ws.Range("C11").Value = "####-10-05"
$ws.Range("D11").Value = "Affitto 4.to Trimestre 2021"
$ws.Range("E11").Value = "external.MISC"
$ws.Range("F11").Value = "z0bug.fy_%(year)s"

# Update the selection to match the new authored state (E3:E11 active)
$ws.Range("E3:E11").Select()
